$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("@prefix")
$ws.Range("A1").Value = "ome"
$ws = $wb.Worksheets.Item("Plate")
$ws.Range("F3").Value = "ome:well"
$ws.Range("G3").Value = "ome:plateAcquisition"
$ws.Range("B4").Value = "ome:Plate"
$ws.Range("F4").Value = "ome:Well"
$ws.Range("G4").Value = "ome:PlateAcquisition"
$ws = $wb.Worksheets.Item("Well")
$ws.Range("D3").Value = "ome:wellSample"
$ws.Range("E3").Value = "ome:row"
$ws.Range("F3").Value = "ome:column"
$ws.Range("G3").Value = "ome:reagent"
$ws.Range("B4").Value = "ome:Well"
$ws.Range("D4").Value = "ome:WellSample"
$ws.Range("G4").Value = "ome:Reagent"
$ws = $wb.Worksheets.Item("Well_Sample")
$ws.Range("D3").Value = "ome:image"
$ws.Range("E3").Value = "ome:index"
$ws.Range("B4").Value = "ome:WellSample"
$ws.Range("D4").Value = "ome:Image"
$ws = $wb.Worksheets.Item("Plate_Acquisiotion")
$ws.Range("D3").Value = "ome:startTime"
$ws.Range("E3").Value = "ome:endTime"
$ws.Range("F3").Value = "ome:wellSample"
$ws.Range("B4").Value = "ome:PlateAcquisition"
$ws.Range("F4").Value = "ome:WellSample"
$ws = $wb.Worksheets.Item("Screen")
$ws.Range("E3").Value = "ome:plate"
$ws.Range("F3").Value = "ome:reagentSet"
$ws.Range("B4").Value = "ome:Screen"
$ws.Range("E4").Value = "ome:Plate"
$ws.Range("F4").Value = "ome:ReagentSet"
$ws = $wb.Worksheets.Item("Reagent_Set")
$ws.Range("C3").Value = "ome:reagent"
$ws.Range("B4").Value = "ome:ReagentSet"
$ws.Range("C4").Value = "ome:Reagent"
$ws = $wb.Worksheets.Item("Reagent")
$ws.Range("B4").Value = "ome:Reagent"
$ws = $wb.Worksheets.Item("Image")
$ws.Range("E3").Value = "ome:pixels"
$ws.Range("F3").Value = "ome:acquisitionDate"
$ws.Range("B4").Value = "ome:Image"
$ws.Range("E4").Value = "ome:Pixels"
$ws = $wb.Worksheets.Item("Pixels")
$ws.Range("D3").Value = "ome:pixelType"
$ws.Range("E3").Value = "ome:dimensionOrder"
$ws.Range("F3").Value = "ome:physicalSizeX"
$ws.Range("G3").Value = "ome:physicalSizeY"
$ws.Range("H3").Value = "ome:sizeC"
$ws.Range("I3").Value = "ome:sizeT"
$ws.Range("J3").Value = "ome:sizeX"
$ws.Range("K3").Value = "ome:sizeY"
$ws.Range("L3").Value = "ome:sizeZ"
$ws.Range("M3").Value = "ome:channel"
$ws.Range("N3").Value = "ome:binData"
$ws.Range("B4").Value = "ome:Pixels"
$ws.Range("D4").Value = "ome:PixelType"
$ws.Range("E4").Value = "ome:DimensionOrder"
$ws.Range("M4").Value = "ome:Channel"
$ws.Range("N4").Value = "ome:BinData"
$ws = $wb.Worksheets.Item("Channel")
$ws.Range("D3").Value = "ome:color"
$ws.Range("B4").Value = "ome:Channel"
$ws.Range("D4").Value = "ome:Color"
$ws = $wb.Worksheets.Item("Color")
$ws.Range("B4").Value = "ome:Color"
$ws = $wb.Worksheets.Item("Binary_Data")
$ws.Range("C3").Value = "ome:bigEndian"
$ws.Range("D3").Value = "ome:data"
$ws.Range("E3").Value = "ome:length"
$ws.Range("B4").Value = "ome:BinData"
